# Generate Report for Handback
#
# This mirrors the "handback" pass of the localization-status report:
# the status string flips from "Ready for handoff" to "Handed back: in
# sync with en-US", the per-language sheets get their "Latest Target
# File" / "Latest Handback File" / "Latest Handback DateTime" columns
# populated (they were previously empty placeholders), and the new
# target-file cells become hyperlinks to the source .md docs (matching
# the existing hyperlinks already present in column A).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1630fdd67acca85fa3b3744ab29ce380270640b6/e2e/3e03fae2-c1e2-4a8b-bc2f-483c0054a3f9.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1630fdd67acca85fa3b3744ab29ce380270640b6/e2e/903d541f-2ac7-4eea-87c8-b8ab02e5adcd.md"
$mdName1 = "3e03fae2-c1e2-4a8b-bc2f-483c0054a3f9.md"
$mdName2 = "903d541f-2ac7-4eea-87c8-b8ab02e5adcd.md"

# ---------------------------------------------------------------------
# Overview sheet: Status columns (E = zh-cn, F = de-de) for both rows
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = 29.9777050018311
$overview.Columns.Item(6).ColumnWidth = 29.9777050018311

# ---------------------------------------------------------------------
# zh-cn / de-de sheets share the same column layout:
#   C = Status, I = Latest Target File, J = Latest Handback File,
#   K = Latest Handback DateTime
# ---------------------------------------------------------------------
$langSheets = @(
    @{
        Name = "zh-cn"
        Xlf1 = "3e03fae2-c1e2-4a8b-bc2f-483c0054a3f9.788205e2bd0e3e5ffd2ec0869b34f23a00e4dc21.zh-cn.xlf"
        Xlf2 = "903d541f-2ac7-4eea-87c8-b8ab02e5adcd.1220c86478874b8b3417a92cfb10294e45c6eef6.zh-cn.xlf"
        HandbackDateTime = "2016-10-19 11:35:56"
    },
    @{
        Name = "de-de"
        Xlf1 = "3e03fae2-c1e2-4a8b-bc2f-483c0054a3f9.788205e2bd0e3e5ffd2ec0869b34f23a00e4dc21.de-de.xlf"
        Xlf2 = "903d541f-2ac7-4eea-87c8-b8ab02e5adcd.1220c86478874b8b3417a92cfb10294e45c6eef6.de-de.xlf"
        HandbackDateTime = "2016-10-19 11:36:15"
    }
)

foreach ($cfg in $langSheets) {
    $sheet = $wb.Worksheets.Item($cfg.Name)

    $sheet.Range("C2").Value = $statusText
    $sheet.Range("C3").Value = $statusText

    # Row 2 -> 3e03fae2... file
    $sheet.Hyperlinks.Add($sheet.Range("I2"), $mdUrl1, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdName1) | Out-Null
    $sheet.Range("I2").Font.Underline = 2
    $sheet.Range("I2").Font.Color = 15570276
    $sheet.Range("J2").Value = $cfg.Xlf1
    $sheet.Range("K2").Value = $cfg.HandbackDateTime

    # Row 3 -> 903d541f... file
    $sheet.Hyperlinks.Add($sheet.Range("I3"), $mdUrl2, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdName2) | Out-Null
    $sheet.Range("I3").Font.Underline = 2
    $sheet.Range("I3").Font.Color = 15570276
    $sheet.Range("J3").Value = $cfg.Xlf2
    $sheet.Range("K3").Value = $cfg.HandbackDateTime

    $sheet.Columns.Item(3).ColumnWidth = 29.9777050018311
    $sheet.Columns.Item(9).ColumnWidth = 40
    $sheet.Columns.Item(10).ColumnWidth = 40
}

Write-Output "Handback report generated."
